$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers, values are IPSA closing levels;
# weekend days use the "--" placeholder already used elsewhere in column B).
$newRows = @(
    @{ Row = 642; Date = 44106; Value = 3666.03 },
    @{ Row = 643; Date = 44107; Value = "--" },
    @{ Row = 644; Date = 44108; Value = "--" },
    @{ Row = 645; Date = 44109; Value = 3675.83 },
    @{ Row = 646; Date = 44110; Value = 3621.71 },
    @{ Row = 647; Date = 44111; Value = 3588.35 },
    @{ Row = 648; Date = 44112; Value = 3650.27 },
    @{ Row = 649; Date = 44113; Value = 3677.45 }
)

foreach ($item in $newRows) {
    $row = $item.Row

    # Copy the date cell format from the row above (A641 carries the date style)
    $ws.Range("A641").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)
    $ws.Range("A" + $row).Value = $item.Date

    if ($item.Value -eq "--") {
        # Reuse the text-placeholder formatting/style from an existing "--" cell (B636)
        $ws.Range("B636").Copy()
        $ws.Range("B" + $row).PasteSpecial(-4122)
        $ws.Range("B" + $row).Value = "--"
    } else {
        # Reuse the numeric value formatting/style from the row above (B641)
        $ws.Range("B641").Copy()
        $ws.Range("B" + $row).PasteSpecial(-4122)
        $ws.Range("B" + $row).Value = $item.Value
    }
}

$excel.CutCopyMode = 0

# Update the named range "IPSA" to cover the new extent
$wb.Names.Item("IPSA").RefersTo = "=IPSA!`$A`$1:`$B`$649"

# Update sheet view: scroll so row 644 is near top, and select B652 (next empty cell)
$excel.ActiveWindow.ScrollRow = 644
$ws.Range("B652").Select()
